$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 2.111846333333333
$ws.Cells.Item(2, 8).Value = 6.335539
$ws.Cells.Item(2, 9).Value = 0.01909882549924913
$ws.Cells.Item(2, 10).Value = 0.01909882549924913
$ws.Cells.Item(2, 13).Value = 172.8922906666667
$ws.Cells.Item(2, 14).Value = 518.676872
$ws.Cells.Item(2, 15).Value = 0.3392537931509703
$ws.Cells.Item(2, 16).Value = 0.3392537931509703
$ws.Cells.Item(2, 17).Value = 365.1219501060009
$ws.Cells.Item(2, 18).Value = 3286.097550954008
$ws.Cells.Item(2, 19).Value = 0.006479348995348742
$ws.Cells.Item(2, 20).Value = 0.006479348995348742

# Row 3
$ws.Cells.Item(3, 7).Value = 2.111846333333333
$ws.Cells.Item(3, 8).Value = 6.335539
$ws.Cells.Item(3, 9).Value = 0.01909882549924913
$ws.Cells.Item(3, 10).Value = 0.01909882549924913
$ws.Cells.Item(3, 15).Value = 0.3999079859779965
$ws.Cells.Item(3, 16).Value = 0.3999079859779965
$ws.Cells.Item(3, 17).Value = 430.4010350100097
$ws.Cells.Item(3, 18).Value = 3873.609315090088
$ws.Cells.Item(3, 19).Value = 0.007637772839949926
$ws.Cells.Item(3, 20).Value = 0.007637772839949926

# Row 4
$ws.Cells.Item(4, 7).Value = 2.111846333333333
$ws.Cells.Item(4, 8).Value = 6.335539
$ws.Cells.Item(4, 9).Value = 0.01909882549924913
$ws.Cells.Item(4, 10).Value = 0.01909882549924913
$ws.Cells.Item(4, 13).Value = 99.55997466666668
$ws.Cells.Item(4, 14).Value = 298.679924
$ws.Cells.Item(4, 15).Value = 0.1953591968817216
$ws.Cells.Item(4, 16).Value = 0.1953591968817216
$ws.Cells.Item(4, 17).Value = 210.2553674465596
$ws.Cells.Item(4, 18).Value = 1892.298307019036
$ws.Cells.Item(4, 19).Value = 0.003731131210917457
$ws.Cells.Item(4, 20).Value = 0.003731131210917457

# Row 5
$ws.Cells.Item(5, 7).Value = 2.111846333333333
$ws.Cells.Item(5, 8).Value = 6.335539
$ws.Cells.Item(5, 9).Value = 0.01909882549924913
$ws.Cells.Item(5, 10).Value = 0.01909882549924913
$ws.Cells.Item(5, 13).Value = 33.36976233333333
$ws.Cells.Item(5, 14).Value = 100.109287
$ws.Cells.Item(5, 15).Value = 0.06547902398931163
$ws.Cells.Item(5, 16).Value = 0.06547902398931162
$ws.Cells.Item(5, 17).Value = 70.47181022785477
$ws.Cells.Item(5, 18).Value = 634.246292050693
$ws.Cells.Item(5, 19).Value = 0.001250572453033011
$ws.Cells.Item(5, 20).Value = 0.001250572453033011

# Row 6
$ws.Cells.Item(6, 7).Value = 51.92481233333333
$ws.Cells.Item(6, 9).Value = 0.4695904783329055
$ws.Cells.Item(6, 10).Value = 0.4695904783329055
$ws.Cells.Item(6, 13).Value = 172.8922906666667
$ws.Cells.Item(6, 14).Value = 518.676872
$ws.Cells.Item(6, 15).Value = 0.3392537931509703
$ws.Cells.Item(6, 16).Value = 0.3392537931509703
$ws.Cells.Item(6, 17).Value = 8977.399746746783
$ws.Cells.Item(6, 18).Value = 80796.59772072105
$ws.Cells.Item(6, 19).Value = 0.1593103510020167
$ws.Cells.Item(6, 20).Value = 0.1593103510020167

# Row 7
$ws.Cells.Item(7, 7).Value = 51.92481233333333
$ws.Cells.Item(7, 9).Value = 0.4695904783329055
$ws.Cells.Item(7, 10).Value = 0.4695904783329055
$ws.Cells.Item(7, 15).Value = 0.3999079859779965
$ws.Cells.Item(7, 16).Value = 0.3999079859779965
$ws.Cells.Item(7, 18).Value = 95241.98497019969
$ws.Cells.Item(7, 19).Value = 0.1877929824245562
$ws.Cells.Item(7, 20).Value = 0.1877929824245562

# Row 8
$ws.Cells.Item(8, 7).Value = 51.92481233333333
$ws.Cells.Item(8, 9).Value = 0.4695904783329055
$ws.Cells.Item(8, 10).Value = 0.4695904783329055
$ws.Cells.Item(8, 13).Value = 99.55997466666668
$ws.Cells.Item(8, 14).Value = 298.679924
$ws.Cells.Item(8, 15).Value = 0.1953591968817216
$ws.Cells.Item(8, 16).Value = 0.1953591968817216
$ws.Cells.Item(8, 17).Value = 5169.633000478088
$ws.Cells.Item(8, 18).Value = 46526.69700430278
$ws.Cells.Item(8, 19).Value = 0.09173881871041992
$ws.Cells.Item(8, 20).Value = 0.09173881871041992

# Row 9
$ws.Cells.Item(9, 7).Value = 51.92481233333333
$ws.Cells.Item(9, 9).Value = 0.4695904783329055
$ws.Cells.Item(9, 10).Value = 0.4695904783329055
$ws.Cells.Item(9, 13).Value = 33.36976233333333
$ws.Cells.Item(9, 14).Value = 100.109287
$ws.Cells.Item(9, 15).Value = 0.06547902398931163
$ws.Cells.Item(9, 16).Value = 0.06547902398931162
$ws.Cells.Item(9, 17).Value = 1732.718646766269
$ws.Cells.Item(9, 18).Value = 15594.46782089642
$ws.Cells.Item(9, 19).Value = 0.03074832619591264
$ws.Cells.Item(9, 20).Value = 0.03074832619591263

# Row 10
$ws.Cells.Item(10, 7).Value = 56.36634066666667
$ws.Cells.Item(10, 8).Value = 169.099022
$ws.Cells.Item(10, 9).Value = 0.5097581615820991
$ws.Cells.Item(10, 10).Value = 0.5097581615820991
$ws.Cells.Item(10, 13).Value = 172.8922906666667
$ws.Cells.Item(10, 14).Value = 518.676872
$ws.Cells.Item(10, 15).Value = 0.3392537931509703
$ws.Cells.Item(10, 16).Value = 0.3392537931509703
$ws.Cells.Item(10, 17).Value = 9745.305754357687
$ws.Cells.Item(10, 18).Value = 87707.75178921918
$ws.Cells.Item(10, 19).Value = 0.1729373899063923
$ws.Cells.Item(10, 20).Value = 0.1729373899063923

# Row 11
$ws.Cells.Item(11, 7).Value = 56.36634066666667
$ws.Cells.Item(11, 8).Value = 169.099022
$ws.Cells.Item(11, 9).Value = 0.5097581615820991
$ws.Cells.Item(11, 10).Value = 0.5097581615820991
$ws.Cells.Item(11, 15).Value = 0.3999079859779965
$ws.Cells.Item(11, 16).Value = 0.3999079859779965
$ws.Cells.Item(11, 17).Value = 11487.64044984656
$ws.Cells.Item(11, 18).Value = 103388.764048619
$ws.Cells.Item(11, 19).Value = 0.2038563597341434
$ws.Cells.Item(11, 20).Value = 0.2038563597341434

# Row 12
$ws.Cells.Item(12, 7).Value = 56.36634066666667
$ws.Cells.Item(12, 8).Value = 169.099022
$ws.Cells.Item(12, 9).Value = 0.5097581615820991
$ws.Cells.Item(12, 10).Value = 0.5097581615820991
$ws.Cells.Item(12, 13).Value = 99.55997466666668
$ws.Cells.Item(12, 14).Value = 298.679924
$ws.Cells.Item(12, 15).Value = 0.1953591968817216
$ws.Cells.Item(12, 16).Value = 0.1953591968817216
$ws.Cells.Item(12, 17).Value = 5611.831448826037
$ws.Cells.Item(12, 18).Value = 50506.48303943433
$ws.Cells.Item(12, 19).Value = 0.09958594505058176
$ws.Cells.Item(12, 20).Value = 0.09958594505058176

# Row 13
$ws.Cells.Item(13, 7).Value = 56.36634066666667
$ws.Cells.Item(13, 8).Value = 169.099022
$ws.Cells.Item(13, 9).Value = 0.5097581615820991
$ws.Cells.Item(13, 10).Value = 0.5097581615820991
$ws.Cells.Item(13, 13).Value = 33.36976233333333
$ws.Cells.Item(13, 14).Value = 100.109287
$ws.Cells.Item(13, 15).Value = 0.06547902398931163
$ws.Cells.Item(13, 16).Value = 0.06547902398931162
$ws.Cells.Item(13, 17).Value = 1880.931391646368
$ws.Cells.Item(13, 18).Value = 16928.38252481731
$ws.Cells.Item(13, 19).Value = 0.03337846689098166
$ws.Cells.Item(13, 20).Value = 0.03337846689098165

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.171671
$ws.Cells.Item(14, 8).Value = 0.5150129999999999
$ws.Cells.Item(14, 9).Value = 0.001552534585746342
$ws.Cells.Item(14, 10).Value = 0.001552534585746342
$ws.Cells.Item(14, 13).Value = 172.8922906666667
$ws.Cells.Item(14, 14).Value = 518.676872
$ws.Cells.Item(14, 15).Value = 0.3392537931509703
$ws.Cells.Item(14, 16).Value = 0.3392537931509703
$ws.Cells.Item(14, 17).Value = 29.68059243103733
$ws.Cells.Item(14, 18).Value = 267.125331879336
$ws.Cells.Item(14, 19).Value = 0.0005267032472125168
$ws.Cells.Item(14, 20).Value = 0.0005267032472125168

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.171671
$ws.Cells.Item(15, 8).Value = 0.5150129999999999
$ws.Cells.Item(15, 9).Value = 0.001552534585746342
$ws.Cells.Item(15, 10).Value = 0.001552534585746342
$ws.Cells.Item(15, 15).Value = 0.3999079859779965
$ws.Cells.Item(15, 16).Value = 0.3999079859779965
$ws.Cells.Item(15, 17).Value = 34.98709868941067
$ws.Cells.Item(15, 18).Value = 314.8838882046959
$ws.Cells.Item(15, 19).Value = 0.0006208709793470028
$ws.Cells.Item(15, 20).Value = 0.0006208709793470028

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.171671
$ws.Cells.Item(16, 8).Value = 0.5150129999999999
$ws.Cells.Item(16, 9).Value = 0.001552534585746342
$ws.Cells.Item(16, 10).Value = 0.001552534585746342
$ws.Cells.Item(16, 13).Value = 99.55997466666668
$ws.Cells.Item(16, 14).Value = 298.679924
$ws.Cells.Item(16, 15).Value = 0.1953591968817216
$ws.Cells.Item(16, 16).Value = 0.1953591968817216
$ws.Cells.Item(16, 17).Value = 17.09156041100134
$ws.Cells.Item(16, 18).Value = 153.824043699012
$ws.Cells.Item(16, 19).Value = 0.0003033019098025017
$ws.Cells.Item(16, 20).Value = 0.0003033019098025017

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.171671
$ws.Cells.Item(17, 8).Value = 0.5150129999999999
$ws.Cells.Item(17, 9).Value = 0.001552534585746342
$ws.Cells.Item(17, 10).Value = 0.001552534585746342
$ws.Cells.Item(17, 13).Value = 33.36976233333333
$ws.Cells.Item(17, 14).Value = 100.109287
$ws.Cells.Item(17, 15).Value = 0.06547902398931163
$ws.Cells.Item(17, 16).Value = 0.06547902398931162
$ws.Cells.Item(17, 17).Value = 5.728620469525667
$ws.Cells.Item(17, 18).Value = 51.55758422573099
$ws.Cells.Item(17, 19).Value = 0.0001016584493843207
$ws.Cells.Item(17, 20).Value = 0.0001016584493843207
